$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19, shifting existing rows down
$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 5
$ws.Cells.Item(19, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(19, 3).Value = "Maule"
$ws.Cells.Item(19, 4).Value = 44670
$ws.Cells.Item(19, 5).Value = 7
$ws.Cells.Item(19, 6).Value = 100112017
$ws.Cells.Item(19, 7).Value = "Apio"
$ws.Cells.Item(19, 8).Value = "Americana (o)"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 600
$ws.Cells.Item(19, 11).Value = 8000
$ws.Cells.Item(19, 12).Value = 8000
$ws.Cells.Item(19, 13).Value = 8000
$ws.Cells.Item(19, 14).Value = "`$/docena de matas"
$ws.Cells.Item(19, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 16).Value = 1333
$ws.Cells.Item(19, 17).Value = 6
$ws.Cells.Item(19, 18).Value = "Hortaliza"
